$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.743.60"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "2.078.71"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.82"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.394"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("D12").Value = "2.384.57"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.81"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.00"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.775"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.36"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "2.114.29"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "37.677.41"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.57"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.60"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -2.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.53"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.137"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.52"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0633"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.52%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.36"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.93"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0974"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.72"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.93%  "
$ws.Range("D45").Value = "1.441.52"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "2.269.40"
$ws.Range("E51").Value = "  -2.18%  "
